$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data has a "project_year" / "country_application_name" table in columns
# A:B (row 1 = headers). Rows are grouped into contiguous blocks per year.
# Sort each year's block of rows by country name (column B) ascending.

$blocks = @(
    @(2, 55),
    @(56, 111),
    @(112, 185),
    @(186, 252),
    @(253, 307)
)

foreach ($block in $blocks) {
    $startRow = $block[0]
    $endRow = $block[1]
    $rng = $ws.Range("A$($startRow):B$($endRow)")
    $keyRng = $ws.Range("B$($startRow):B$($endRow)")
    $rng.Sort($keyRng)
}
